# Applies the "remove Gamelogic project, modify SLG building config" edit.
# This inserts two new columns (Icon, ShowName) before the existing Desc
# column on Sheet1 of the BB_Build workbook, and fills in appropriate
# values for the new columns in every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G currently holds "Desc". Move its contents to column I (without
# doing a structural column-insert, which would also drag every later
# column's width metadata along with it) so that G/H are free for the new
# "Icon"/"ShowName" columns.
for ($row = 1; $row -le 10; $row++) {
    $oldDesc = $ws.Cells.Item($row, 7).Value()
    $ws.Cells.Item($row, 9).Value = $oldDesc
}

# Header row
$ws.Cells.Item(1, 7).Value = "Icon"
$ws.Cells.Item(1, 8).Value = "ShowName"

# Data rows: Icon (G) mirrors the prefab's short name, ShowName (H) mirrors
# the existing localized Desc text that now lives in column I.
$iconValues = @{
    2  = "Altar_1_1"
    3  = "Arena_1_1"
    4  = "Camp_1_1"
    5  = "GoldMine_1_1"
    6  = "Item_hourse_1_1"
    7  = "League_1_1"
    8  = "MagicHourse_1_1"
    9  = "Tower_1_1"
    10 = "Town_1_1"
}

for ($row = 2; $row -le 10; $row++) {
    $descValue = $ws.Cells.Item($row, 9).Value()
    $ws.Cells.Item($row, 7).Value = $iconValues[$row]
    $ws.Cells.Item($row, 8).Value = $descValue

    # These rows use the workbook's "text" cell style (numFmtId 49) on every
    # other column (A, D, E, F, ...), so apply the same formatting to the
    # new G/H columns to match.
    $ws.Cells.Item($row, 7).NumberFormat = "@"
    $ws.Cells.Item($row, 8).NumberFormat = "@"
    $ws.Cells.Item($row, 9).NumberFormat = "@"
}

# Column widths: G, H, I should all share width 11 (matching the target
# layout where the old width-14 "H" spec is dropped). Excel stores column
# widths in units derived from the workbook's "Normal" font maximum digit
# width (7px here), so request 72/7 characters to obtain an exact stored
# width of 11.
$ws.Columns.Item(7).ColumnWidth = 72 / 7
$ws.Columns.Item(8).ColumnWidth = 72 / 7
$ws.Columns.Item(9).ColumnWidth = 72 / 7

# Restore the selection to match the saved workbook state.
$ws.Range("H10").Select() | Out-Null
